$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.961.71"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.364.05"
$ws.Range("E3").Value = "  +2.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "

# Row 7
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  -0.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.57%  "

# Row 12
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0785"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15
$ws.Range("D15").Value = "2.731.98"
$ws.Range("E15").Value = "  +2.26%  "

# Row 16
$ws.Range("D16").Value = "2.350.06"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.73%  "

# Row 18
$ws.Range("D18").Value = "42.944.93"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("E19").Value = "  -2.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.21%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").Value = "  -0.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.93"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

# Row 24
$ws.Range("E24").Value = "  -4.32%  "

# Row 25
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("E26").Value = "  +0.66%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "

# Row 28
$ws.Range("E28").Value = "  +0.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.81%  "

# Row 31
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0718"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.97%  "

# Row 35
$ws.Range("E35").Value = "  +4.23%  "

# Row 36
$ws.Range("E36").Value = "  +3.70%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "127.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -23.55%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.97%  "

# Row 40
$ws.Range("E40").Value = "  +2.48%  "

# Row 41
$ws.Range("E41").Value = "  -0.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.26%  "

# Row 43
$ws.Range("D43").Value = "1.934.06"
$ws.Range("E43").Value = "  +0.47%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0278"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.27%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.11%  "

# Row 47
$ws.Range("E47").Value = "  -0.97%  "

# Row 48
$ws.Range("D48").Value = "2.589.78"
$ws.Range("E48").Value = "  +1.95%  "

# Row 49
$ws.Range("E49").Value = "  +2.17%  "

# Row 50
$ws.Range("E50").Value = "  +1.86%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.21%  "
